# Applies the Denmark Superliga 2023-2024 update:
#  1) Three pairs of adjacent match rows were reordered (rows 21/22, 39/40,
#     57/58) - the "home"..."url" data (columns F:V) for each pair is
#     swapped while the Indice/pais/torneio/temporada/data columns (A:E)
#     stay put.
#  2) Six new match rows were appended (rows 74-79), extending the used
#     range from A1:V73 to A1:V79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $r1, $r2, $startCol, $endCol) {
    $range1 = $ws.Range("$startCol$r1" + ":" + "$endCol$r1")
    $range2 = $ws.Range("$startCol$r2" + ":" + "$endCol$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value = $v2
    $range2.Value = $v1
}

# --- 1) Swap the F:V data between the three row pairs -----------------
Swap-RowData $ws 21 22 "F" "V"
Swap-RowData $ws 39 40 "F" "V"
Swap-RowData $ws 57 58 "F" "V"

# --- 2) Append the six new rows (74-79) --------------------------------
$newRows = @(
    @{ row = 74; vals = @{
        A = 73; B = "denmark"; C = "superliga"; D = "2023-2024"; E = 45226.79166666666
        F = "Midtjylland"; G = 2; H = "Lyngby"; I = 1
        J = 1.57; K = "23/10/2023 19:12"; L = 1.56; M = "27/10/2023 18:49"
        N = 4.22; O = "23/10/2023 19:12"; P = 4.25; Q = "27/10/2023 18:56"
        R = 5.73; S = "23/10/2023 19:12"; T = 6.14; U = "27/10/2023 18:56"
        V = "https://www.betexplorer.com/football/denmark/superliga/midtjylland-lyngby/CWP1koUj/"
    } },
    @{ row = 75; vals = @{
        A = 74; B = "denmark"; C = "superliga"; D = "2023-2024"; E = 45227.70833333334
        F = "FC Copenhagen"; G = 4; H = "Hvidovre IF"; I = 0
        J = 1.19; K = "21/10/2023 18:13"; L = 1.19; M = "28/10/2023 16:56"
        N = 7.33; O = "21/10/2023 18:13"; P = 7.77; Q = "28/10/2023 16:58"
        R = 13.58; S = "21/10/2023 18:13"; T = 13.41; U = "28/10/2023 16:58"
        V = "https://www.betexplorer.com/football/denmark/superliga/fc-copenhagen-hvidovre-if/YuRcjREp/"
    } },
    @{ row = 76; vals = @{
        A = 75; B = "denmark"; C = "superliga"; D = "2023-2024"; E = 45228.58333333334
        F = "Vejle"; G = 1; H = "Viborg"; I = 1
        J = 2.66; K = "22/10/2023 16:12"; L = 2.96; M = "29/10/2023 13:58"
        N = 3.44; O = "22/10/2023 16:12"; P = 3.31; Q = "29/10/2023 13:53"
        R = 2.54; S = "22/10/2023 16:12"; T = 2.55; U = "29/10/2023 13:58"
        V = "https://www.betexplorer.com/football/denmark/superliga/vejle-viborg/tp3U7gpI/"
    } },
    @{ row = 77; vals = @{
        A = 76; B = "denmark"; C = "superliga"; D = "2023-2024"; E = 45228.66666666666
        F = "Silkeborg"; G = 0; H = "Odense"; I = 0
        J = 1.68; K = "23/10/2023 19:12"; L = 1.71; M = "29/10/2023 15:59"
        N = 4.15; O = "23/10/2023 19:12"; P = 4.1; Q = "29/10/2023 15:59"
        R = 4.74; S = "23/10/2023 19:12"; T = 4.82; U = "29/10/2023 15:59"
        V = "https://www.betexplorer.com/football/denmark/superliga/silkeborg-odense/ll9pEU6i/"
    } },
    @{ row = 78; vals = @{
        A = 77; B = "denmark"; C = "superliga"; D = "2023-2024"; E = 45228.75
        F = "Brondby"; G = 2; H = "Nordsjaelland"; I = 1
        J = 2.3; K = "23/10/2023 09:12"; L = 2.69; M = "29/10/2023 17:51"
        N = 3.69; O = "23/10/2023 09:12"; P = 3.53; Q = "29/10/2023 17:56"
        R = 2.81; S = "23/10/2023 09:12"; T = 2.64; U = "29/10/2023 17:51"
        V = "https://www.betexplorer.com/football/denmark/superliga/brondby-nordsjaelland/rapR8ZUB/"
    } },
    @{ row = 79; vals = @{
        A = 78; B = "denmark"; C = "superliga"; D = "2023-2024"; E = 45229.79166666666
        F = "Aarhus"; G = 2; H = "Randers FC"; I = 1
        J = 1.6; K = "23/10/2023 09:12"; L = 1.96; M = "30/10/2023 18:58"
        N = 3.88; O = "23/10/2023 09:12"; P = 3.32; Q = "30/10/2023 18:58"
        R = 5.39; S = "23/10/2023 09:12"; T = 4.48; U = "30/10/2023 18:58"
        V = "https://www.betexplorer.com/football/denmark/superliga/aarhus-randers-fc/0zeM9FF5/"
    } }
)

$lastTemplateRow = 73

foreach ($item in $newRows) {
    $r = $item.row

    # Copy the number-format / border / bold styling used by column A
    # (Indice) and column E (data_partida) from the last existing row so
    # the new rows look like the rest of the table.
    $ws.Range("A$lastTemplateRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("E$lastTemplateRow").Copy() | Out-Null
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null

    foreach ($col in $item.vals.Keys) {
        $addr = "$col$r"
        $ws.Range($addr).Value = $item.vals[$col]
    }
}

$excel.CutCopyMode = $false
